$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C1 currently holds "Wit2". Replace it with the rich-text value "Աբ2", where
# "Աբ" is rendered in the Sylfaen font (matching the font used elsewhere in
# this sheet, e.g. the armenian text cells) and "2" stays in the regular
# Calibri font.
$rng = $ws.Range("C1")
$rng.Value2 = "Աբ2"
$rng.Characters(1, 2).Font.Name = "Sylfaen"
$rng.Characters(3, 1).Font.Name = "Calibri"

# Row 1 grows a bit taller once it carries the mixed-font content (matches
# the height already used by every other row on this sheet).
$ws.Rows.Item(1).RowHeight = 17

# Move the active selection from D5 to C1.
$null = $ws.Range("C1").Select()
